$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.378.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.563.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.17%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.576"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.575.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.35%  "
$ws.Range("E10").Value = "  -5.20%  "
$ws.Range("E11").Value = "  -5.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.60%  "
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.011.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.353.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.566.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.81%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.417"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.674.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.159"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0816"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.06%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.944"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.79%  "
$ws.Range("E38").Value = "  -6.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.853"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "291.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.997"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.606"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0536"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.94%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0228"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.72%  "
